$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = -16.46144747230769
$ws.Cells.Item(2, 3).Value = -16.46144747230769
$ws.Cells.Item(2, 4).Value = -16.46144747230769
$ws.Cells.Item(2, 5).Value = -16.46144747230769
$ws.Cells.Item(2, 6).Value = -16.46144747230769
$ws.Cells.Item(2, 7).Value = -16.46144747230769
$ws.Cells.Item(2, 8).Value = -16.46144747230769
$ws.Cells.Item(2, 9).Value = -16.46144747230769
$ws.Cells.Item(2, 10).Value = -16.46144747230769
$ws.Cells.Item(2, 11).Value = -16.46144747230769
$ws.Cells.Item(3, 2).Value = -16.46144747230769
$ws.Cells.Item(3, 3).Value = -16.46144747230769
$ws.Cells.Item(3, 4).Value = -16.46144747230769
$ws.Cells.Item(3, 5).Value = -16.46144747230769
$ws.Cells.Item(3, 6).Value = -16.46144747230769
$ws.Cells.Item(3, 7).Value = -16.46144747230769
$ws.Cells.Item(3, 8).Value = -16.46144747230769
$ws.Cells.Item(3, 9).Value = 2.730366001868378
$ws.Cells.Item(3, 10).Value = -16.46144747230769
$ws.Cells.Item(3, 11).Value = -16.46144747230769
$ws.Cells.Item(4, 2).Value = -16.46144747230769
$ws.Cells.Item(4, 3).Value = -16.46144747230769
$ws.Cells.Item(4, 4).Value = 2.771921766834931
$ws.Cells.Item(4, 5).Value = -16.46144747230769
$ws.Cells.Item(4, 6).Value = 3.112667417731577
$ws.Cells.Item(4, 7).Value = -16.46144747230769
$ws.Cells.Item(4, 8).Value = 1.708815034269539
$ws.Cells.Item(4, 9).Value = -16.46144747230769
$ws.Cells.Item(4, 10).Value = 2.146536559003594
$ws.Cells.Item(4, 11).Value = -16.46144747230769
$ws.Cells.Item(5, 2).Value = -16.46144747230769
$ws.Cells.Item(5, 3).Value = -16.46144747230769
$ws.Cells.Item(5, 4).Value = -16.46144747230769
$ws.Cells.Item(5, 5).Value = -16.46144747230769
$ws.Cells.Item(5, 6).Value = -16.46144747230769
$ws.Cells.Item(5, 7).Value = 2.859690111365378
$ws.Cells.Item(5, 8).Value = -16.46144747230769
$ws.Cells.Item(5, 9).Value = -16.46144747230769
$ws.Cells.Item(5, 10).Value = -16.46144747230769
$ws.Cells.Item(5, 11).Value = -16.46144747230769
$ws.Cells.Item(6, 2).Value = -16.46144747230769
$ws.Cells.Item(6, 3).Value = -16.46144747230769
$ws.Cells.Item(6, 4).Value = -16.46144747230769
$ws.Cells.Item(6, 5).Value = -16.46144747230769
$ws.Cells.Item(6, 6).Value = -16.46144747230769
$ws.Cells.Item(6, 7).Value = -16.46144747230769
$ws.Cells.Item(6, 8).Value = -16.46144747230769
$ws.Cells.Item(6, 9).Value = -16.46144747230769
$ws.Cells.Item(6, 10).Value = -16.46144747230769
$ws.Cells.Item(6, 11).Value = -16.46144747230769
$ws.Cells.Item(7, 2).Value = 2.509739607330871
$ws.Cells.Item(7, 3).Value = -16.46144747230769
$ws.Cells.Item(7, 4).Value = -16.46144747230769
$ws.Cells.Item(7, 5).Value = -16.46144747230769
$ws.Cells.Item(7, 6).Value = -16.46144747230769
$ws.Cells.Item(7, 7).Value = -16.46144747230769
$ws.Cells.Item(7, 8).Value = -16.46144747230769
$ws.Cells.Item(7, 9).Value = -16.46144747230769
$ws.Cells.Item(7, 10).Value = -16.46144747230769
$ws.Cells.Item(7, 11).Value = -16.46144747230769
$ws.Cells.Item(8, 2).Value = -16.46144747230769
$ws.Cells.Item(8, 3).Value = -16.46144747230769
$ws.Cells.Item(8, 4).Value = -16.46144747230769
$ws.Cells.Item(8, 5).Value = 1.58766188939
$ws.Cells.Item(8, 6).Value = -16.46144747230769
$ws.Cells.Item(8, 7).Value = -16.46144747230769
$ws.Cells.Item(8, 8).Value = -16.46144747230769
$ws.Cells.Item(8, 9).Value = -16.46144747230769
$ws.Cells.Item(8, 10).Value = -16.46144747230769
$ws.Cells.Item(8, 11).Value = -16.46144747230769
$ws.Cells.Item(9, 2).Value = 3.838410050940509
$ws.Cells.Item(9, 3).Value = -16.46144747230769
$ws.Cells.Item(9, 4).Value = -16.46144747230769
$ws.Cells.Item(9, 5).Value = -16.46144747230769
$ws.Cells.Item(9, 6).Value = -16.46144747230769
$ws.Cells.Item(9, 7).Value = -16.46144747230769
$ws.Cells.Item(9, 8).Value = -16.46144747230769
$ws.Cells.Item(9, 9).Value = -16.46144747230769
$ws.Cells.Item(9, 10).Value = -16.46144747230769
$ws.Cells.Item(9, 11).Value = -16.46144747230769
$ws.Cells.Item(10, 2).Value = -16.46144747230769
$ws.Cells.Item(10, 3).Value = -16.46144747230769
$ws.Cells.Item(10, 4).Value = -16.46144747230769
$ws.Cells.Item(10, 5).Value = -16.46144747230769
$ws.Cells.Item(10, 6).Value = -16.46144747230769
$ws.Cells.Item(10, 7).Value = -16.46144747230769
$ws.Cells.Item(10, 8).Value = -16.46144747230769
$ws.Cells.Item(10, 9).Value = 1.31272829900185
$ws.Cells.Item(10, 10).Value = -16.46144747230769
$ws.Cells.Item(10, 11).Value = 1.963444607981134
$ws.Cells.Item(11, 2).Value = -16.46144747230769
$ws.Cells.Item(11, 3).Value = -16.46144747230769
$ws.Cells.Item(11, 4).Value = -16.46144747230769
$ws.Cells.Item(11, 5).Value = 3.159690430022193
$ws.Cells.Item(11, 6).Value = -16.46144747230769
$ws.Cells.Item(11, 7).Value = 2.794588407602761
$ws.Cells.Item(11, 8).Value = -16.46144747230769
$ws.Cells.Item(11, 9).Value = -16.46144747230769
$ws.Cells.Item(11, 10).Value = -16.46144747230769
$ws.Cells.Item(11, 11).Value = 1.812043102464373
$ws.Cells.Item(12, 2).Value = -16.46144747230769
$ws.Cells.Item(12, 3).Value = -16.46144747230769
$ws.Cells.Item(12, 4).Value = -16.46144747230769
$ws.Cells.Item(12, 5).Value = -16.46144747230769
$ws.Cells.Item(12, 6).Value = -16.46144747230769
$ws.Cells.Item(12, 7).Value = -16.46144747230769
$ws.Cells.Item(12, 8).Value = -16.46144747230769
$ws.Cells.Item(12, 9).Value = -16.46144747230769
$ws.Cells.Item(12, 10).Value = -16.46144747230769
$ws.Cells.Item(12, 11).Value = -16.46144747230769
$ws.Cells.Item(13, 2).Value = -16.46144747230769
$ws.Cells.Item(13, 3).Value = -16.46144747230769
$ws.Cells.Item(13, 4).Value = -16.46144747230769
$ws.Cells.Item(13, 5).Value = 2.323305573916441
$ws.Cells.Item(13, 6).Value = -16.46144747230769
$ws.Cells.Item(13, 7).Value = -16.46144747230769
$ws.Cells.Item(13, 8).Value = -16.46144747230769
$ws.Cells.Item(13, 9).Value = -16.46144747230769
$ws.Cells.Item(13, 10).Value = 2.103658325867963
$ws.Cells.Item(13, 11).Value = 1.910673642395263
$ws.Cells.Item(14, 2).Value = -16.46144747230769
$ws.Cells.Item(14, 3).Value = -16.46144747230769
$ws.Cells.Item(14, 4).Value = 1.773748728774624
$ws.Cells.Item(14, 5).Value = -16.46144747230769
$ws.Cells.Item(14, 6).Value = -16.46144747230769
$ws.Cells.Item(14, 7).Value = -16.46144747230769
$ws.Cells.Item(14, 8).Value = -16.46144747230769
$ws.Cells.Item(14, 9).Value = -16.46144747230769
$ws.Cells.Item(14, 10).Value = -16.46144747230769
$ws.Cells.Item(14, 11).Value = 2.156253143648072
$ws.Cells.Item(15, 2).Value = -16.46144747230769
$ws.Cells.Item(15, 3).Value = -16.46144747230769
$ws.Cells.Item(15, 4).Value = 1.632800351920362
$ws.Cells.Item(15, 5).Value = -16.46144747230769
$ws.Cells.Item(15, 6).Value = -16.46144747230769
$ws.Cells.Item(15, 7).Value = -16.46144747230769
$ws.Cells.Item(15, 8).Value = -16.46144747230769
$ws.Cells.Item(15, 9).Value = -16.46144747230769
$ws.Cells.Item(15, 10).Value = -16.46144747230769
$ws.Cells.Item(15, 11).Value = -16.46144747230769
$ws.Cells.Item(16, 2).Value = -16.46144747230769
$ws.Cells.Item(16, 3).Value = -16.46144747230769
$ws.Cells.Item(16, 4).Value = -16.46144747230769
$ws.Cells.Item(16, 5).Value = -16.46144747230769
$ws.Cells.Item(16, 6).Value = -16.46144747230769
$ws.Cells.Item(16, 7).Value = -16.46144747230769
$ws.Cells.Item(16, 8).Value = -16.46144747230769
$ws.Cells.Item(16, 9).Value = -16.46144747230769
$ws.Cells.Item(16, 10).Value = 2.199081223893042
$ws.Cells.Item(16, 11).Value = -16.46144747230769
$ws.Cells.Item(17, 2).Value = -16.46144747230769
$ws.Cells.Item(17, 3).Value = -16.46144747230769
$ws.Cells.Item(17, 4).Value = 1.710867712093027
$ws.Cells.Item(17, 5).Value = -16.46144747230769
$ws.Cells.Item(17, 6).Value = -16.46144747230769
$ws.Cells.Item(17, 7).Value = -16.46144747230769
$ws.Cells.Item(17, 8).Value = 1.372486466871122
$ws.Cells.Item(17, 9).Value = 1.845170540284616
$ws.Cells.Item(17, 10).Value = 1.918645735765221
$ws.Cells.Item(17, 11).Value = -16.46144747230769
$ws.Cells.Item(18, 2).Value = -16.46144747230769
$ws.Cells.Item(18, 3).Value = -16.46144747230769
$ws.Cells.Item(18, 4).Value = -16.46144747230769
$ws.Cells.Item(18, 5).Value = -16.46144747230769
$ws.Cells.Item(18, 6).Value = -16.46144747230769
$ws.Cells.Item(18, 7).Value = -16.46144747230769
$ws.Cells.Item(18, 8).Value = 1.582973330172867
$ws.Cells.Item(18, 9).Value = 1.318765468036413
$ws.Cells.Item(18, 10).Value = 1.536884019448322
$ws.Cells.Item(18, 11).Value = -16.46144747230769
$ws.Cells.Item(19, 2).Value = -16.46144747230769
$ws.Cells.Item(19, 3).Value = -16.46144747230769
$ws.Cells.Item(19, 4).Value = 0.9302744072515426
$ws.Cells.Item(19, 5).Value = -16.46144747230769
$ws.Cells.Item(19, 6).Value = -16.46144747230769
$ws.Cells.Item(19, 7).Value = -16.46144747230769
$ws.Cells.Item(19, 8).Value = 1.651343257659814
$ws.Cells.Item(19, 9).Value = 1.493595055175653
$ws.Cells.Item(19, 10).Value = -16.46144747230769
$ws.Cells.Item(19, 11).Value = -16.46144747230769
$ws.Cells.Item(20, 2).Value = -16.46144747230769
$ws.Cells.Item(20, 3).Value = -16.46144747230769
$ws.Cells.Item(20, 4).Value = 0.5557197908064442
$ws.Cells.Item(20, 5).Value = -16.46144747230769
$ws.Cells.Item(20, 6).Value = 3.504619969822
$ws.Cells.Item(20, 7).Value = -16.46144747230769
$ws.Cells.Item(20, 8).Value = 1.889340672647954
$ws.Cells.Item(20, 9).Value = 0.9826696143003706
$ws.Cells.Item(20, 10).Value = -16.46144747230769
$ws.Cells.Item(20, 11).Value = 2.128002521954439
$ws.Cells.Item(21, 2).Value = -16.46144747230769
$ws.Cells.Item(21, 3).Value = 4.32191290652433
$ws.Cells.Item(21, 4).Value = -16.46144747230769
$ws.Cells.Item(21, 5).Value = 1.610250557558669
$ws.Cells.Item(21, 6).Value = -16.46144747230769
$ws.Cells.Item(21, 7).Value = 2.536771063852566
$ws.Cells.Item(21, 8).Value = 2.104327736438449
$ws.Cells.Item(21, 9).Value = -16.46144747230769
$ws.Cells.Item(21, 10).Value = -16.46144747230769
$ws.Cells.Item(21, 11).Value = -16.46144747230769
